# Auto-generated script to update C2:C13 values across all 30 SBFL metric sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tarantula")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.09477124183006536
$ws.Range("C4").Value = 0.1977124183006536
$ws.Range("C5").Value = 0.2173202614379085
$ws.Range("C6").Value = 0.2532679738562091
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.2450980392156863
$ws.Range("C9").Value = 0.4509803921568628
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5392156862745098
$ws.Range("C12").Value = 7.980392156862745
$ws.Range("C13").Value = 129.9607843137255

$ws = $wb.Worksheets.Item("Ochiai")
$ws.Range("C2").Value = 0.09967320261437909
$ws.Range("C3").Value = 0.1650326797385622
$ws.Range("C4").Value = 0.2385620915032678
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2499999999999998
$ws.Range("C7").Value = 0.2549019607843137
$ws.Range("C8").Value = 0.4215686274509804
$ws.Range("C9").Value = 0.6078431372549019
$ws.Range("C10").Value = 0.6078431372549019
$ws.Range("C11").Value = 0.6274509803921569
$ws.Range("C12").Value = 4.84313725490196
$ws.Range("C13").Value = 121.6960784313726

$ws = $wb.Worksheets.Item("Op2")
$ws.Range("C2").Value = 0.008169934640522875
$ws.Range("C3").Value = 0.03431372549019608
$ws.Range("C4").Value = 0.04901960784313725
$ws.Range("C5").Value = 0.04901960784313725
$ws.Range("C6").Value = 0.05392156862745097
$ws.Range("C7").Value = 0.0196078431372549
$ws.Range("C8").Value = 0.07843137254901961
$ws.Range("C9").Value = 0.107843137254902
$ws.Range("C10").Value = 0.107843137254902
$ws.Range("C11").Value = 0.1176470588235294
$ws.Range("C12").Value = 15.93137254901961
$ws.Range("C13").Value = 130.1470588235294

$ws = $wb.Worksheets.Item("Barinel")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.09477124183006536
$ws.Range("C4").Value = 0.1977124183006536
$ws.Range("C5").Value = 0.2173202614379085
$ws.Range("C6").Value = 0.2532679738562091
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.2450980392156863
$ws.Range("C9").Value = 0.4509803921568628
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5392156862745098
$ws.Range("C12").Value = 7.980392156862745
$ws.Range("C13").Value = 129.9607843137255

$ws = $wb.Worksheets.Item("Dstar")
$ws.Range("C2").Value = 0.09967320261437909
$ws.Range("C3").Value = 0.1650326797385622
$ws.Range("C4").Value = 0.2385620915032678
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2581699346405227
$ws.Range("C7").Value = 0.2549019607843137
$ws.Range("C8").Value = 0.4215686274509804
$ws.Range("C9").Value = 0.6078431372549019
$ws.Range("C10").Value = 0.6078431372549019
$ws.Range("C11").Value = 0.6274509803921569
$ws.Range("C12").Value = 4.784313725490196
$ws.Range("C13").Value = 122

$ws = $wb.Worksheets.Item("Russell_rao")
$ws.Range("C12").Value = 22.58823529411765
$ws.Range("C13").Value = 144.7450980392157

$ws = $wb.Worksheets.Item("Simple_matching")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Rogers_tanimoto")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Ample")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1535947712418302
$ws.Range("C4").Value = 0.2467320261437907
$ws.Range("C5").Value = 0.2630718954248363
$ws.Range("C6").Value = 0.2843137254901957
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.392156862745098
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6372549019607843
$ws.Range("C11").Value = 0.6862745098039216
$ws.Range("C12").Value = 4.46078431372549
$ws.Range("C13").Value = 152.8529411764706

$ws = $wb.Worksheets.Item("Jaccard")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.901960784313726
$ws.Range("C13").Value = 118.5490196078431

$ws = $wb.Worksheets.Item("Cohen")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.882352941176471
$ws.Range("C13").Value = 225.4803921568627

$ws = $wb.Worksheets.Item("Scott")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2009803921568627
$ws.Range("C5").Value = 0.2009803921568627
$ws.Range("C6").Value = 0.2009803921568627
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5098039215686274
$ws.Range("C10").Value = 0.5098039215686274
$ws.Range("C11").Value = 0.5098039215686274
$ws.Range("C12").Value = 33.6078431372549
$ws.Range("C13").Value = 323.8333333333333

$ws = $wb.Worksheets.Item("Rogot1")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2009803921568627
$ws.Range("C5").Value = 0.2009803921568627
$ws.Range("C6").Value = 0.2009803921568627
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5098039215686274
$ws.Range("C10").Value = 0.5098039215686274
$ws.Range("C11").Value = 0.5098039215686274
$ws.Range("C12").Value = 33.6078431372549
$ws.Range("C13").Value = 323.8333333333333

$ws = $wb.Worksheets.Item("Geometric_mean")
$ws.Range("C2").Value = 0.09967320261437909
$ws.Range("C3").Value = 0.1650326797385622
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2352941176470586
$ws.Range("C6").Value = 0.2467320261437906
$ws.Range("C7").Value = 0.2549019607843137
$ws.Range("C8").Value = 0.4215686274509804
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.5980392156862745
$ws.Range("C11").Value = 0.6274509803921569
$ws.Range("C12").Value = 5.225490196078431
$ws.Range("C13").Value = 221.5294117647059

$ws = $wb.Worksheets.Item("M2")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1568627450980393
$ws.Range("C4").Value = 0.2434640522875815
$ws.Range("C5").Value = 0.2728758169934637
$ws.Range("C6").Value = 0.2941176470588231
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.4019607843137255
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6666666666666666
$ws.Range("C11").Value = 0.7156862745098039
$ws.Range("C12").Value = 4.794117647058823
$ws.Range("C13").Value = 128.7843137254902

$ws = $wb.Worksheets.Item("Wong1")
$ws.Range("C12").Value = 22.58823529411765
$ws.Range("C13").Value = 144.7450980392157

$ws = $wb.Worksheets.Item("Sokal")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Sorensen_dice")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.901960784313726
$ws.Range("C13").Value = 118.5490196078431

$ws = $wb.Worksheets.Item("Dice")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.901960784313726
$ws.Range("C13").Value = 118.5490196078431

$ws = $wb.Worksheets.Item("Humman")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Wong2")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Euclid")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Zoltar")
$ws.Range("C2").Value = 0.02777777777777778
$ws.Range("C3").Value = 0.065359477124183
$ws.Range("C4").Value = 0.0915032679738562
$ws.Range("C5").Value = 0.09967320261437909
$ws.Range("C6").Value = 0.09967320261437909
$ws.Range("C7").Value = 0.07843137254901961
$ws.Range("C8").Value = 0.1764705882352941
$ws.Range("C9").Value = 0.2450980392156863
$ws.Range("C10").Value = 0.2647058823529412
$ws.Range("C11").Value = 0.2647058823529412
$ws.Range("C12").Value = 10.33333333333333
$ws.Range("C13").Value = 125.8725490196078

$ws = $wb.Worksheets.Item("Rogot2")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1535947712418302
$ws.Range("C4").Value = 0.2467320261437907
$ws.Range("C5").Value = 0.2630718954248363
$ws.Range("C6").Value = 0.2843137254901957
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.392156862745098
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6372549019607843
$ws.Range("C11").Value = 0.6862745098039216
$ws.Range("C12").Value = 4.411764705882353
$ws.Range("C13").Value = 214.843137254902

$ws = $wb.Worksheets.Item("Hamming")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2173202614379085
$ws.Range("C5").Value = 0.2385620915032678
$ws.Range("C6").Value = 0.2598039215686274
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 0.5
$ws.Range("C12").Value = 96.45098039215686
$ws.Range("C13").Value = 359.5686274509804

$ws = $wb.Worksheets.Item("Fleiss")
$ws.Range("C2").Value = 0.06535947712418298
$ws.Range("C3").Value = 0.1307189542483661
$ws.Range("C4").Value = 0.2124183006535947
$ws.Range("C5").Value = 0.2205882352941175
$ws.Range("C6").Value = 0.2254901960784312
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.3333333333333333
$ws.Range("C9").Value = 0.5392156862745098
$ws.Range("C10").Value = 0.5588235294117647
$ws.Range("C11").Value = 0.5686274509803921
$ws.Range("C12").Value = 31.56862745098039
$ws.Range("C13").Value = 325.1568627450981

$ws = $wb.Worksheets.Item("Anderberg")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.901960784313726
$ws.Range("C13").Value = 118.5490196078431

$ws = $wb.Worksheets.Item("Goodman")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1421568627450981
$ws.Range("C4").Value = 0.2352941176470586
$ws.Range("C5").Value = 0.2434640522875815
$ws.Range("C6").Value = 0.2549019607843134
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.3627450980392157
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6176470588235294
$ws.Range("C11").Value = 0.6470588235294118
$ws.Range("C12").Value = 6.901960784313726
$ws.Range("C13").Value = 118.5490196078431

$ws = $wb.Worksheets.Item("Harmonic_mean")
$ws.Range("C2").Value = 0.07679738562091501
$ws.Range("C3").Value = 0.1535947712418302
$ws.Range("C4").Value = 0.2467320261437907
$ws.Range("C5").Value = 0.2630718954248363
$ws.Range("C6").Value = 0.2843137254901957
$ws.Range("C7").Value = 0.196078431372549
$ws.Range("C8").Value = 0.392156862745098
$ws.Range("C9").Value = 0.5980392156862745
$ws.Range("C10").Value = 0.6372549019607843
$ws.Range("C11").Value = 0.6862745098039216
$ws.Range("C12").Value = 4.411764705882353
$ws.Range("C13").Value = 221.5980392156863

$ws = $wb.Worksheets.Item("Kulczynski2")
$ws.Range("C2").Value = 0.06862745098039214
$ws.Range("C3").Value = 0.1552287581699348
$ws.Range("C4").Value = 0.2254901960784312
$ws.Range("C5").Value = 0.2254901960784312
$ws.Range("C6").Value = 0.2303921568627449
$ws.Range("C7").Value = 0.1764705882352941
$ws.Range("C8").Value = 0.392156862745098
$ws.Range("C9").Value = 0.5686274509803921
$ws.Range("C10").Value = 0.5686274509803921
$ws.Range("C11").Value = 0.5784313725490197
$ws.Range("C12").Value = 7.235294117647059
$ws.Range("C13").Value = 132.3529411764706

